$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road
$ws.Range("H17").Value = 207056
$ws.Range("J17").Value = 213702.97
$ws.Range("L17").Value = 641108.91
$ws.Range("N17").Value = -641444.91
# Row 18: You Grow, Girl
$ws.Range("H18").Value = 1093.25
$ws.Range("I18").Value = 1093.25
$ws.Range("K18").Value = 1093.25
$ws.Range("M18").Value = -809.25
# Row 33: Glazed and Confused
$ws.Range("H33").Value = 3965843.2
$ws.Range("I33").Value = 5451601
$ws.Range("K33").Value = 5451601
$ws.Range("M33").Value = -5451372
# Row 86: Filling in the Blanks
$ws.Range("I86").Value = 166670750
$ws.Range("J86").Value = 41670970
$ws.Range("K86").Value = 166670750
$ws.Range("L86").Value = 41670970
$ws.Range("M86").Value = -166669627
$ws.Range("N86").Value = -41673216
# Row 89: Ink into Antiquity (L)
$ws.Range("I89").Value = 166670750
$ws.Range("J89").Value = 41670970
$ws.Range("K89").Value = 833353750
$ws.Range("L89").Value = 208354850
$ws.Range("M89").Value = -833348134
$ws.Range("N89").Value = -208366082
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1286.2041
$ws.Range("I132").Value = 1350.5946
$ws.Range("K132").Value = 4051.7838
$ws.Range("M132").Value = -1521.7838
# Row 138: All-night Crafting
$ws.Range("H138").Value = 4737.491
$ws.Range("I138").Value = 6092.433
$ws.Range("J138").Value = 3232
$ws.Range("K138").Value = 18277.299
$ws.Range("L138").Value = 9696
$ws.Range("M138").Value = -13137.299
$ws.Range("N138").Value = -19976
# Row 141: Remedy for Reason
$ws.Range("H141").Value = 1966.8422
$ws.Range("I141").Value = 1398.125
$ws.Range("K141").Value = 4194.375
$ws.Range("M141").Value = 985.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 24: A Firm Hand
$ws.Range("H24").Value = 78451.664
$ws.Range("J24").Value = 78451.664
$ws.Range("L24").Value = 78451.664
$ws.Range("N24").Value = -79199.664
# Row 54: Family Secrets
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 329043.56
$ws.Range("I74").Value = 1628.7307
$ws.Range("J74").Value = 1225126.4
$ws.Range("K74").Value = 1628.7307
$ws.Range("L74").Value = 1225126.4
$ws.Range("M74").Value = -754.7307000000001
$ws.Range("N74").Value = -1226874.4
# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 329043.56
$ws.Range("I77").Value = 1628.7307
$ws.Range("J77").Value = 1225126.4
$ws.Range("K77").Value = 8143.6535
$ws.Range("L77").Value = 6125632
$ws.Range("M77").Value = -3775.6535
$ws.Range("N77").Value = -6134368
# Row 100: En Garde and on Guard
$ws.Range("H100").Value = 78451.664
$ws.Range("J100").Value = 78451.664
$ws.Range("L100").Value = 78451.664
$ws.Range("N100").Value = -80615.664
# Row 107: Shielding the Realm
$ws.Range("H107").Value = 39997
$ws.Range("J107").Value = 39997
$ws.Range("L107").Value = 39997
$ws.Range("N107").Value = -47677
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 1277.8182
$ws.Range("I122").Value = 867.375
$ws.Range("K122").Value = 2602.125
$ws.Range("M122").Value = -152.125
# Row 123: The Armoire Is Open
$ws.Range("H123").Value = 99499.5
$ws.Range("J123").Value = 99499.5
$ws.Range("L123").Value = 99499.5
$ws.Range("N123").Value = -109299.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107: The Gold Experience
$ws.Range("H107").Value = 19420.35
$ws.Range("I107").Value = 20179.316
$ws.Range("K107").Value = 20179.316
$ws.Range("M107").Value = -18259.316

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22: Driving Up the Wall
$ws.Range("H22").Value = 1415.25
$ws.Range("I22").Value = 1415.25
$ws.Range("K22").Value = 1415.25
$ws.Range("M22").Value = -1065.25
# Row 31: Wall Not Found
$ws.Range("H31").Value = 2455.3845
$ws.Range("I31").Value = 2381.087
$ws.Range("K31").Value = 2381.087
$ws.Range("M31").Value = -2086.087
# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 2455.3845
$ws.Range("I34").Value = 2381.087
$ws.Range("K34").Value = 2381.087
$ws.Range("M34").Value = -2179.087
# Row 50: The Arsenal of Theocracy
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 60: Bowing to Greater Power
$ws.Range("H60").Value = 11999
$ws.Range("I60").Value = 11999
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 11999
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = -11488
$ws.Range("N60").ClearContents()
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 2105.5
$ws.Range("I134").Value = 2025.4615
$ws.Range("K134").Value = 6076.3845
$ws.Range("M134").Value = -3541.3845

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 57: The Egg Files
$ws.Range("H57").Value = 10416.667
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
# Row 113: Can't Eat Just One
$ws.Range("H113").Value = 1659.8
$ws.Range("J113").Value = 2433
$ws.Range("L113").Value = 7299
$ws.Range("N113").Value = -11639
# Row 117: A Good Omen
$ws.Range("H117").Value = 166669710
$ws.Range("J117").Value = 166669710
$ws.Range("L117").Value = 500009130
$ws.Range("N117").Value = -500016014
# Row 129: Comfort Food
$ws.Range("H129").Value = 13377.25
$ws.Range("I129").Value = 1399.5
$ws.Range("J129").Value = 19366.125
$ws.Range("K129").Value = 4198.5
$ws.Range("L129").Value = 58098.375
$ws.Range("M129").Value = 801.5
$ws.Range("N129").Value = -68098.375
# Row 138: Bring Me Your Tacos
$ws.Range("H138").Value = 3614.5
$ws.Range("I138").Value = 3614.5
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 10843.5
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -5703.5
$ws.Range("N138").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Range("H2").Value = 102.25
$ws.Range("I2").Value = 91.57143000000001
$ws.Range("J2").Value = 117.2
$ws.Range("K2").Value = 91.57143000000001
$ws.Range("L2").Value = 117.2
$ws.Range("M2").Value = 21.42856999999999
$ws.Range("N2").Value = -343.2
# Row 49: Faith and Fashion
$ws.Range("H49").Value = 42889.45
$ws.Range("I49").Value = 37200
$ws.Range("J49").Value = 44785.934
$ws.Range("K49").Value = 37200
$ws.Range("L49").Value = 44785.934
$ws.Range("M49").Value = -37016
$ws.Range("N49").Value = -45153.934
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 8499
$ws.Range("I113").Value = 8000
$ws.Range("J113").Value = 8748.5
$ws.Range("K113").Value = 8000
$ws.Range("L113").Value = 8748.5
$ws.Range("M113").Value = -5830
$ws.Range("N113").Value = -13088.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs
$ws.Range("H22").Value = 3540.2593
$ws.Range("I22").Value = 1231.1428
$ws.Range("J22").Value = 4348.45
$ws.Range("K22").Value = 1231.1428
$ws.Range("L22").Value = 4348.45
$ws.Range("M22").Value = -936.1428000000001
$ws.Range("N22").Value = -4938.45
# Row 27: Fire and Hide
$ws.Range("H27").Value = 3540.2593
$ws.Range("I27").Value = 1231.1428
$ws.Range("J27").Value = 4348.45
$ws.Range("K27").Value = 1231.1428
$ws.Range("L27").Value = 4348.45
$ws.Range("M27").Value = -1124.1428
$ws.Range("N27").Value = -4562.45

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display
$ws.Range("H96").Value = 3668.0625
$ws.Range("I96").Value = 3653.6365
$ws.Range("K96").Value = 3653.6365
$ws.Range("M96").Value = -2280.6365
# Row 113: A Tender Table
$ws.Range("H113").Value = 623.0714
$ws.Range("J113").Value = 753.8
$ws.Range("L113").Value = 2261.4
$ws.Range("N113").Value = -6601.4
# Row 117: The Hunt Continues
$ws.Range("H117").Value = 74000
$ws.Range("J117").Value = 74000
$ws.Range("L117").Value = 74000
$ws.Range("N117").Value = -83178
